$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.984.77'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '2.355.29'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +4.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.03'
$ws.Range('E6').Value = '  +3.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.53'
$ws.Range('E7').Value = '  +11.50%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +18.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('E10').Value = '  +5.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.43'
$ws.Range('E11').Value = '  +10.76%  '
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('D13').Value = '2.706.05'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('E14').Value = '  +10.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.69'
$ws.Range('E15').Value = '  +7.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.909'
$ws.Range('E16').Value = '  +8.56%  '
$ws.Range('D17').Value = '2.356.68'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '43.932.36'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('E19').Value = '  +5.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '77.89'
$ws.Range('E20').Value = '  +5.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.44'
$ws.Range('E21').Value = '  +4.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '254.74'
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  -3.78%  '
$ws.Range('E25').Value = '  +3.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.48'
$ws.Range('E26').Value = '  +6.76%  '
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.59'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('E30').Value = '  +7.61%  '
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('E32').Value = '  +5.41%  '
$ws.Range('E33').Value = '  +3.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0727'
$ws.Range('E34').Value = '  +6.56%  '
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('E36').Value = '  +10.56%  '
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('E39').Value = '  +7.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.47'
$ws.Range('E40').Value = '  +9.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.84'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('E44').Value = '  +4.57%  '
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.36'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('E47').Value = '  +2.06%  '
$ws.Range('E48').Value = '  +12.52%  '
$ws.Range('E49').Value = '  +4.63%  '
$ws.Range('D50').Value = '1.437.94'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('E51').Value = '  +1.41%  '
